$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Make room: insert two new rows above the first data row so the
#    existing 5 rows of data shift from rows 2-6 down to rows 4-8.
# ------------------------------------------------------------------
$ws.Range("A2:C3").Insert(-4121)  # xlShiftDown
$ws.Range("A2:C3").Font.Bold = $false

# ------------------------------------------------------------------
# 2) New rows 2-3 (new shapefiles discovered before the existing ones)
# ------------------------------------------------------------------
$ws.Range("A2").Value = 43850
$ws.Range("B2").Value = "states_il_mo_01"
$ws.Range("C2").Value = "Shapefile for the states of Illinois and Missouri combined."

$ws.Range("A3").Value = 43850
$ws.Range("B3").Value = "stl_msa"
$ws.Range("C3").Value = "Shapefile for the St. Louis Metropolitan Statistical Area (MSA)."

$ws.Range("A2:A3").NumberFormat = "dd\-mmm\-yyyy"
$ws.Range("A2:C3").HorizontalAlignment = -4131  # xlLeft
$ws.Range("A2:C3").WrapText = $true

# ------------------------------------------------------------------
# 3) Rows 4-8 keep the same filenames but get updated notes text
#    (the data that used to live in rows 2-6 before the insert).
# ------------------------------------------------------------------
$ws.Range("C4").Value = "Shapefile for St. Louis MSA with census tracts."
$ws.Range("C5").Value = "Shapefile for St. Louis MSA with census tracts and social data part 01."
$ws.Range("C6").Value = "Shapefile for St. Louis MSA with census tracts and social data part 02."
$ws.Range("C7").Value = "Shapefile for St. Louis MSA with census tracts removed and social data and indicators added."
$ws.Range("C8").Value = "Shapefile for St. Louis MSA with census tracts removed, social and indicators added, and projected to UTM-15."

# ------------------------------------------------------------------
# 4) Append the new rows of data (9-22) for the homework assignment
#    and class project downloads.
# ------------------------------------------------------------------
$newRows = @(
    @(43867, "mc_blk",     "Shapefile for St. Louis MSA mean center of Black population."),
    @(43867, "mc_lat",     "Shapefile for St. Louis MSA mean center of Latino and Hispanic population."),
    @(43867, "mc_wht",     "Shapefile for St. Louis MSA mean center of White population."),
    @(43867, "sd_blk",     "Shapefile for St. Louis MSA standard distance of Black population."),
    @(43867, "sd_lat",     "Shapefile for St. Louis MSA standard distance of Latino and Hispanic population."),
    @(43867, "sd_wht",     "Shapefile for St. Louis MSA standard distance of White population."),
    @(43867, "sde_blk",    "Shapefile for St. Louis MSA directional distribution of Black population."),
    @(43867, "sde_lat",    "Shapefile for St. Louis MSA directional distribution of Latino and Hispanic population."),
    @(43867, "sde_wht",    "Shapefile for St. Louis MSA directional distribution of White population."),
    @(43874, "stl_city",   "Shapefile for the boundary of the City of St. Louis, Missouri."),
    @(43874, "stl_grid01", "Shapefile for the City of St. Louis, Missouri with grid overlay."),
    @(43874, "stl_grid02", "Shapefile for the City of St. Louis, Missouri with grid clipped to city boundary."),
    @(43874, "stl_int01",  "Shapefile for the City of St. Louis, Missouri with grid clipped to city boundary, partial grids removed, and census tracts added."),
    @(43874, "stl_int02",  "Shapefile for the City of St. Louis, Missouri with grid clipped to city boundary, partial grids removed, census tracts added, and population interpolated to grid.")
)

$r = 9
foreach ($row in $newRows) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $r++
}

# Taller rows for the two longest wrapped notes (21 & 22)
$ws.Rows("21:22").RowHeight = 30

# ------------------------------------------------------------------
# 5) Update the frozen-pane view / active selection to match the new
#    bottom of the list.
# ------------------------------------------------------------------
[void]$ws.Range("A5").Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$ws.Range("A23").Select()

Write-Host "edit complete"
